# further cleaning to metadata
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample protocol id shared string (G2:G41): E7760 -> E7420
$ws.Range("G2:G41").Value = "E7420"

# Give the G column cells their own distinct font (Arial 11, black) so they
# pick up a new style slot instead of sharing the previous one.
$gRange = $ws.Range("G2:G41")
$gRange.Font.Name = "Arial"
$gRange.Font.Size = 11
$gRange.Font.Color = 0

# Column H (roboticS2Prep) becomes an explicit =FALSE() formula instead of a
# bare boolean literal.
$ws.Range("H2:H41").Formula = "=FALSE()"

# Restore the view/selection that was captured with the edit.
[void]$ws.Range("G2:G41").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
